$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 41 (Leve Item ID 5478)
$ws.Range("H41").Value2 = 101.833336
$ws.Range("I41").Value2 = 95
$ws.Range("J41").Value2 = 105.25
$ws.Range("K41").Value2 = 95
$ws.Range("L41").Value2 = 105.25
$ws.Range("M41").Value2 = 345
$ws.Range("N41").Value2 = -985.25

# Row 55 (Leve Item ID 5517)
$ws.Range("H55").Value2 = 83334300
$ws.Range("I55").Value2 = 125001304
$ws.Range("J55").Value2 = 300
$ws.Range("K55").Value2 = 125001304
$ws.Range("L55").Value2 = 300
$ws.Range("M55").Value2 = -125001090
$ws.Range("N55").Value2 = -728

# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value2 = 2963.5417
$ws.Range("I116").Value2 = 2918.0557
$ws.Range("J116").Value2 = 3100
$ws.Range("K116").Value2 = 2918.0557
$ws.Range("L116").Value2 = 3100
$ws.Range("M116").Value2 = 523.9443000000001
$ws.Range("N116").Value2 = -9984

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value2 = 2965.4546
$ws.Range("I132").Value2 = 2917.4138
$ws.Range("J132").Value2 = 3313.75
$ws.Range("K132").Value2 = 8752.241399999999
$ws.Range("L132").Value2 = 9941.25
$ws.Range("M132").Value2 = -6222.241399999999
$ws.Range("N132").Value2 = -15001.25

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value2 = 3790.0715
$ws.Range("I137").Value2 = 3447.8635
$ws.Range("J137").Value2 = 5044.8335
$ws.Range("K137").Value2 = 10343.5905
$ws.Range("L137").Value2 = 15134.5005
$ws.Range("M137").Value2 = -7793.5905
$ws.Range("N137").Value2 = -20234.5005

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value2 = 192753.89
$ws.Range("I138").Value2 = 2357.3684
$ws.Range("K138").Value2 = 7072.1052
$ws.Range("M138").Value2 = -1932.1052

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value2 = 479496.88
$ws.Range("I32").Value2 = 629466.25
$ws.Range("K32").Value2 = 629466.25
$ws.Range("M32").Value2 = -629179.25

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value2 = 1068.1666
$ws.Range("I94").Value2 = 969.6667
$ws.Range("K94").Value2 = 969.6667
$ws.Range("M94").Value2 = -518.6667

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value2 = 3456.6155
$ws.Range("I134").Value2 = 3215.111
$ws.Range("K134").Value2 = 9645.332999999999
$ws.Range("M134").Value2 = -7110.332999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 21 (Leve Item ID 2000)
$ws.Range("H21").Value2 = 501
$ws.Range("I21").Value2 = 501
$ws.Range("K21").Value2 = 501
$ws.Range("M21").Value2 = -266

# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value2 = 2058.7273
$ws.Range("I31").Value2 = 1075.3846
$ws.Range("J31").Value2 = 3479.111
$ws.Range("K31").Value2 = 1075.3846
$ws.Range("L31").Value2 = 3479.111
$ws.Range("M31").Value2 = -780.3846000000001
$ws.Range("N31").Value2 = -4069.111

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value2 = 2058.7273
$ws.Range("I34").Value2 = 1075.3846
$ws.Range("J34").Value2 = 3479.111
$ws.Range("K34").Value2 = 1075.3846
$ws.Range("L34").Value2 = 3479.111
$ws.Range("M34").Value2 = -873.3846000000001
$ws.Range("N34").Value2 = -3883.111

# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value2 = 1296.8572
$ws.Range("I58").Value2 = 920
$ws.Range("K58").Value2 = 920
$ws.Range("M58").Value2 = -717

# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value2 = 508.2
$ws.Range("I107").Value2 = 417.11765
$ws.Range("J107").Value2 = 627.3077
$ws.Range("K107").Value2 = 417.11765
$ws.Range("L107").Value2 = 627.3077
$ws.Range("M107").Value2 = 1502.88235
$ws.Range("N107").Value2 = -4467.3077

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value2 = 6411695.5
$ws.Range("I132").Value2 = 805.4706
$ws.Range("J132").Value2 = 18521154
$ws.Range("K132").Value2 = 2416.4118
$ws.Range("L132").Value2 = 55563462
$ws.Range("M132").Value2 = 113.5882000000001
$ws.Range("N132").Value2 = -55568522

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value2 = 1716.7333
$ws.Range("I134").Value2 = 1519.3077
$ws.Range("K134").Value2 = 4557.9231
$ws.Range("M134").Value2 = -2022.9231

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value2 = 1296.8572
$ws.Range("I136").Value2 = 920
$ws.Range("K136").Value2 = 2760
$ws.Range("M136").Value2 = -210

$ws = $wb.Worksheets.Item("CUL")
# Row 108 (Leve Item ID 27853)
$ws.Range("H108").Value2 = 599.8
$ws.Range("I108").Value2 = 599.8
$ws.Range("J108").Value2 = 0
$ws.Range("K108").Value2 = 1799.4
$ws.Range("L108").Value2 = 0
$ws.Range("M108").ClearContents()
$ws.Range("N108").Value2 = 1080.6

# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value2 = 3440.883
$ws.Range("I132").Value2 = 2540.1482
$ws.Range("K132").Value2 = 22861.3338
$ws.Range("M132").Value2 = -20331.3338

$ws = $wb.Worksheets.Item("GSM")
# Row 19 (Leve Item ID 2668)
$ws.Range("H19").Value2 = 27501.166
$ws.Range("I19").Value2 = 5001
$ws.Range("J19").Value2 = 32001.2
$ws.Range("K19").Value2 = 5001
$ws.Range("L19").Value2 = 32001.2
$ws.Range("M19").Value2 = -4713
$ws.Range("N19").Value2 = -32577.2

# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value2 = 2702
$ws.Range("I126").Value2 = 2878
$ws.Range("J126").Value2 = 2584.6667
$ws.Range("K126").Value2 = 8634
$ws.Range("L126").Value2 = 7754.000100000001
$ws.Range("M126").Value2 = -6164
$ws.Range("N126").Value2 = -12694.0001

# Row 135 (Leve Item ID 42006)
$ws.Range("H135").Value2 = 80000
$ws.Range("J135").Value2 = 80000
$ws.Range("L135").Value2 = 80000
$ws.Range("N135").Value2 = -90140

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (Leve Item ID 5289)
$ws.Range("H16").Value2 = 3942.8572
$ws.Range("I16").Value2 = 3942.8572
$ws.Range("J16").Value2 = 0
$ws.Range("K16").Value2 = 3942.8572
$ws.Range("L16").Value2 = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value2 = -3772.8572

# Row 61 (Leve Item ID 27740)
$ws.Range("H61").Value2 = 5011.533
$ws.Range("I61").Value2 = 5231.8
$ws.Range("K61").Value2 = 5231.8
$ws.Range("M61").Value2 = -5029.8

# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value2 = 7596.0527
$ws.Range("I93").Value2 = 10326.167
$ws.Range("J93").Value2 = 2915.8572
$ws.Range("K93").Value2 = 10326.167
$ws.Range("L93").Value2 = 2915.8572
$ws.Range("M93").Value2 = -9078.166999999999
$ws.Range("N93").Value2 = -5411.8572

# Row 113 (Leve Item ID 27740)
$ws.Range("H113").Value2 = 5011.533
$ws.Range("I113").Value2 = 5231.8
$ws.Range("K113").Value2 = 5231.8
$ws.Range("M113").Value2 = -3061.8

# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value2 = 2313.125
$ws.Range("I122").Value2 = 2300
$ws.Range("J122").Value2 = 2321
$ws.Range("K122").Value2 = 6900
$ws.Range("L122").Value2 = 6963
$ws.Range("M122").Value2 = -4450
$ws.Range("N122").Value2 = -11863

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (Leve Item ID 12596)
$ws.Range("H81").Value2 = 7978.1665
$ws.Range("J81").Value2 = 6499.5
$ws.Range("L81").Value2 = 12999
$ws.Range("N81").Value2 = -15121

# Row 84 (Leve Item ID 12596)
$ws.Range("H84").Value2 = 7978.1665
$ws.Range("J84").Value2 = 6499.5
$ws.Range("L84").Value2 = 64995
$ws.Range("N84").Value2 = -75603

# Row 96 (Leve Item ID 19977)
$ws.Range("H96").Value2 = 3960
$ws.Range("I96").Value2 = 3097.1428
$ws.Range("K96").Value2 = 3097.1428
$ws.Range("M96").Value2 = -1724.1428

# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value2 = 1100.65
$ws.Range("I113").Value2 = 1294.2
$ws.Range("J113").Value2 = 520
$ws.Range("K113").Value2 = 3882.6
$ws.Range("L113").Value2 = 1560
$ws.Range("M113").Value2 = -1712.6
$ws.Range("N113").Value2 = -5900

# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value2 = 1095.4857
$ws.Range("I126").Value2 = 717.26086
$ws.Range("J126").Value2 = 1820.4166
$ws.Range("K126").Value2 = 2151.78258
$ws.Range("L126").Value2 = 5461.2498
$ws.Range("M126").Value2 = 318.2174199999999
$ws.Range("N126").Value2 = -10401.2498

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value2 = 5211189
$ws.Range("I132").Value2 = 3012.2856
$ws.Range("K132").Value2 = 9036.856800000001
$ws.Range("M132").Value2 = -6506.856800000001
